$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2 is a text cell ("41" -> "38"). Force text format so Excel doesn't
# auto-convert the numeric-looking string into a real number, then reset
# the cell style back to Normal so no stray number-format style sticks.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "38"
$ws.Range("A2").Style = "Normal"

$ws.Range("C2").Value = "svmC:14.2543794902gamma:4.0"

$ws.Range("D2").Value = 0.6920353982300885
$ws.Range("E2").Value = 0.8442087814170319
$ws.Range("F2").Value = 0.476991150442478
$ws.Range("G2").Value = 0.476991150442478
$ws.Range("H2").Value = 0.9070796460176991
$ws.Range("I2").Value = 0.6540357286002186
$ws.Range("J2").Value = 0.6030891999362946
$ws.Range("K2").Value = 0.6030891999362946
$ws.Range("L2").Value = 0.4286918432896029
$ws.Range("M2").Value = 0.7379716500900617
$ws.Range("N2").Value = 539
$ws.Range("O2").Value = 591
$ws.Range("P2").Value = 105
$ws.Range("Q2").Value = 1025
